# Add two new "Software Requirements" slides (Functional / Non-functional)
# right after the existing "Server keeps a cached list..." slide, and merge
# two runs on the "Server Issues" slide into one.

$p = $ppt.ActivePresentation

# --- Step 1: create the two new slides by duplicating slide 11 ("Server
# Implementation" / cached-list slide), which already has the bare <a:bodyPr/>
# (no autofit override) on both its body and title placeholders, matching the
# target slides' formatting most closely.
$src = $p.Slides.Item(11)

$reqFunctional = $src.Duplicate().Item(1)
$reqNonFunctional = $src.Duplicate().Item(1)

# Put them in presentation order: ... slide9, Functional, Non-functional, old10, old11, ...
$reqFunctional.MoveTo(10)
$reqNonFunctional.MoveTo(11)

# --- Step 2: fill in the "Functional" requirements slide content ---
$bodyA = $reqFunctional.Shapes.Item(1).TextFrame.TextRange
$bodyA.Text = "Functional:`rServer should be robust and be resilient to failure`rServer log messages should be clear on activity`rDevice crashes should not corrupt any part of server`rDevices may only have one " + [char]8220 + "owner" + [char]8221 + " at any given instance`rClients can actively control only one device any given instance`rServer must be secure against unwarranted input"

for ($i = 2; $i -le $bodyA.Paragraphs().Count; $i++) {
    $bodyA.Paragraphs($i).IndentLevel = 2
}

$reqFunctional.Shapes.Item(2).TextFrame.TextRange.Text = "Software Requirements"

# --- Step 3: fill in the "Non-functional" requirements slide content ---
$bodyB = $reqNonFunctional.Shapes.Item(1).TextFrame.TextRange
$bodyB.Text = "Non-functional:`rReal-time devices require near real-time feedback`rDevices should have minimal setup to boot up and connect to server`rSystem should be responsive under any amount of stress`rServer deployment should be straight forward`r`r"

for ($i = 2; $i -le $bodyB.Paragraphs().Count; $i++) {
    $bodyB.Paragraphs($i).IndentLevel = 2
}

$reqNonFunctional.Shapes.Item(2).TextFrame.TextRange.Text = "Software Requirements"

# --- Step 4: merge the two runs on the "Server Issues" slide into one run ---
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $titleShape = $null
    foreach ($sh in $slide.Shapes) {
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText -and $sh.TextFrame.TextRange.Text -eq "Server Issues") {
            $titleShape = $sh
        }
    }
    if ($titleShape -ne $null) {
        foreach ($sh in $slide.Shapes) {
            if ($sh -ne $titleShape -and $sh.HasTextFrame) {
                $tr = $sh.TextFrame.TextRange
                if ($tr.Text.StartsWith("Expanding and extending")) {
                    $tr.Paragraphs(1).Text = "Expanding and extending the functionality of the API can be difficult, server can store what available commands exist for a certain device, but the user interface cannot dynamically translate commands into a practical layout `r"
                }
            }
        }
    }
}
